$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1. Mark column N ("Purchased?") with "y" for rows 10-42.
#    Most rows naturally pick up the correct style (2, from the column
#    default) just by setting the value. A subset of rows need style 1
#    (matching the "Normal 2" cell style used by column A on row 10),
#    which we copy over explicitly with PasteSpecial (formats only) so
#    we reuse the existing style index instead of creating a new one.
# ----------------------------------------------------------------------
for ($r = 10; $r -le 42; $r++) {
    $ws.Cells.Item($r, 14).Value = "y"
}

$style1Rows = @(14, 15, 18, 19, 22, 23, 26, 27, 28, 29, 32, 33, 36, 37, 40, 41)
$ws.Range("A10").Copy()
foreach ($r in $style1Rows) {
    $ws.Cells.Item($r, 14).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# 2. Add the LCSC / JLC purchase summary rows below the existing Ebay
#    row (56), copying the number format from B56 for the price cells.
# ----------------------------------------------------------------------
$ws.Range("A57").Value = "JLC"
$ws.Range("B57").Value = 7.49

$ws.Range("A58").Value = "LCSC"
$ws.Range("B58").Value = 12.18

$ws.Range("A60").Value = "Total"
$ws.Range("B60").Formula = "=SUM(B56:B58)"

$ws.Range("B56").Copy()
$ws.Range("B57:B58").PasteSpecial(-4122)
$ws.Range("B60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A56").Copy()
$ws.Range("A57:A58").PasteSpecial(-4122)
$ws.Range("A60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# 3. Cosmetic view changes recorded in the diff (zoom level & selection).
# ----------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("O44").Select() | Out-Null
